$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "23.872.30"
Set-TextValue "E2" "  -0.29%  "
Set-TextValue "D3" "1.647.58"
Set-TextValue "E3" "  +1.42%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.35%  "
Set-TextValue "D5" "308.76"
Set-TextValue "E5" "  +0.13%  "
Set-TextValue "E6" "  +0.45%  "
Set-TextValue "E7" "  -1.03%  "
Set-TextValue "E8" "  -0.75%  "
Set-TextValue "D9" "51.36"
Set-TextValue "E9" "  +2.45%  "
Set-TextValue "E10" "  -1.43%  "
Set-TextValue "D11" "1.002"
Set-TextValue "E11" "  +0.39%  "
Set-TextValue "D12" "0.08418"
Set-TextValue "E12" "  -0.84%  "
Set-TextValue "D13" "23.84"
Set-TextValue "E13" "  -0.43%  "
Set-TextValue "D14" "7.061"
Set-TextValue "E14" "  -0.56%  "
Set-TextValue "D15" "7.889"
Set-TextValue "E15" "  +3.18%  "
Set-TextValue "D16" "0.00001312"
Set-TextValue "E16" "  +1.78%  "
Set-TextValue "D17" "1.649.04"
Set-TextValue "E17" "  +1.82%  "
Set-TextValue "D18" "94.35"
Set-TextValue "D19" "0.06976"
Set-TextValue "E19" "  +0.74%  "
Set-TextValue "E20" "  -2.76%  "
Set-TextValue "D21" "6.920"
Set-TextValue "E21" "  +0.68%  "
Set-TextValue "E22" "  +0.40%  "
Set-TextValue "D23" "13.65"
Set-TextValue "E23" "  +1.24%  "
Set-TextValue "D24" "23.863.02"
Set-TextValue "E24" "  -0.42%  "
Set-TextValue "D25" "2.446"
Set-TextValue "E25" "  -0.95%  "
Set-TextValue "D26" "2.950"
Set-TextValue "E26" "  +2.54%  "
Set-TextValue "D27" "21.95"
Set-TextValue "E27" "  -1.49%  "
Set-TextValue "D28" "150.78"
Set-TextValue "E28" "  -3.89%  "
Set-TextValue "D29" "5.394"
Set-TextValue "E29" "  +1.65%  "
Set-TextValue "D30" "138.25"
Set-TextValue "E30" "  -1.87%  "
Set-TextValue "D31" "7.764"
Set-TextValue "E31" "  -2.24%  "
Set-TextValue "D32" "2.515"
Set-TextValue "E32" "  +1.41%  "
Set-TextValue "D33" "1.832.16"
Set-TextValue "E33" "  +1.82%  "
Set-TextValue "D34" "1.041"
Set-TextValue "E34" "  +4.51%  "
Set-TextValue "D35" "0.08018"
Set-TextValue "E35" "  -1.71%  "
Set-TextValue "D36" "0.02942"
Set-TextValue "E36" "  +0.83%  "
Set-TextValue "D37" "6.681"
Set-TextValue "E37" "  +0.46%  "
Set-TextValue "E38" "  +3.71%  "
Set-TextValue "D39" "0.2670"
Set-TextValue "E39" "  -0.31%  "
Set-TextValue "D40" "0.09084"
Set-TextValue "E40" "  -0.93%  "
Set-TextValue "D41" "0.7540"
Set-TextValue "E41" "  -0.54%  "
Set-TextValue "D42" "13.41"
Set-TextValue "E42" "  -2.34%  "
Set-TextValue "E43" "  -1.01%  "
Set-TextValue "D44" "16.28"
Set-TextValue "E44" "  +1.21%  "
Set-TextValue "D45" "0.6936"
Set-TextValue "E45" "  -0.27%  "
Set-TextValue "D46" "2.450"
Set-TextValue "E46" "  -1.28%  "
Set-TextValue "D47" "4.080"
Set-TextValue "E47" "  +0.05%  "
Set-TextValue "E48" "  +0.43%  "
Set-TextValue "D49" "0.08273"
Set-TextValue "E49" "  -0.37%  "
Set-TextValue "D50" "133.91"
Set-TextValue "E50" "  -1.81%  "
Set-TextValue "D51" "1.205"
Set-TextValue "E51" "  -0.74%  "
